#
# "Convert all v1 to v2" -- reshape the d3fend library workbook from the
# old 2-sheet "v1" layout (library_content, controls) into the new
# 3-sheet "v2" layout (library_meta, controls_meta, controls_content).
#

$wb = $excel.ActiveWorkbook

# --- 0. Grab the original sheets before we start renaming them -------------
$libSheet = $wb.Worksheets.Item(1)      # "library_content" -> "library_meta"
$ctrlSheet = $wb.Worksheets.Item(2)     # "controls"        -> "controls_meta"

# --- 1. Create the new "controls_content" sheet right after the controls
#        sheet, then copy all of the existing reference-control rows into
#        it before we wipe/repurpose the controls sheet itself. ------------
$contentSheet = $wb.Worksheets.Add($null, $ctrlSheet)
$contentSheet.Name = "controls_content"

$ctrlSheet.UsedRange.Copy($contentSheet.Range("A1"))

# --- 2. Rename the two original sheets to their v2 names -------------------
$libSheet.Name = "library_meta"
$ctrlSheet.Name = "controls_meta"

# --- 3. Rebuild "library_meta": drop the "library_" prefix on every key,
#        insert a new leading "type"/"library" row, and drop the old
#        "tab" / "reference_control_base_urn" rows (they move to the new
#        controls_meta sheet). -----------------------------------------------
$copyrightText = @'
Terms of Use
LICENSE
The MITRE Corporation (MITRE) hereby grants you a non-exclusive, royalty-free license to use D3FEND for research, development, and commercial purposes. Any copy you make for such purposes is authorized provided that you reproduce MITRE’s copyright designation and this license in any such copy.
DISCLAIMERS
ALL DOCUMENTS AND THE INFORMATION CONTAINED THEREIN ARE PROVIDED ON AN "AS IS" BASIS AND THE CONTRIBUTOR, THE ORGANIZATION HE/SHE REPRESENTS OR IS SPONSORED BY (IF ANY), THE MITRE CORPORATION, ITS BOARD OF TRUSTEES, OFFICERS, AGENTS, AND EMPLOYEES, DISCLAIM ALL WARRANTIES, EXPRESS OR IMPLIED, INCLUDING BUT NOT LIMITED TO ANY WARRANTY THAT THE USE OF THE INFORMATION THEREIN WILL NOT INFRINGE ANY RIGHTS OR ANY IMPLIED WARRANTIES OF MERCHANTABILITY OR FITNESS FOR A PARTICULAR PURPOSE.
'@

$descriptionText = @'
A cybersecurity ontology designed to standardize vocabulary for employing techniques to counter malicious cyber threats.
Version - 1.0.0 - 2024-12-20
https://d3fend.mitre.org/resources/
'@

$libSheet.Cells.Clear()

$libData = @(
    @("type", "library"),
    @("urn", "urn:intuitem:risk:library:mitre-d3fend"),
    @("version", "1"),
    @("locale", "en"),
    @("publication_date", "2025-01-22"),
    @("ref_id", "d3fend"),
    @("name", "Mitre D3FEND"),
    @("description", $descriptionText),
    @("copyright", $copyrightText),
    @("provider", "Mitre D3FEND"),
    @("packager", "intuitem")
)

# Rows whose B value must stay plain text (would otherwise be auto-coerced
# into a number/date by COM's smart cell-value parsing).
$textValueRows = @(3, 5)

for ($i = 0; $i -lt $libData.Length; $i++) {
    $row = $i + 1
    $libSheet.Cells.Item($row, 1).Value = $libData[$i][0]
    if ($textValueRows -contains $row) {
        $libSheet.Cells.Item($row, 2).NumberFormat = "@"
    }
    $libSheet.Cells.Item($row, 2).Value = $libData[$i][1]
}

# --- 4. Rebuild "controls_meta" with the small reference-controls header
#        that used to live in the "tab" / "reference_control_base_urn" rows
#        of the library sheet. -----------------------------------------------
$ctrlSheet.Cells.Clear()

$ctrlSheet.Cells.Item(1, 1).Value = "type"
$ctrlSheet.Cells.Item(1, 2).Value = "reference_controls"
$ctrlSheet.Cells.Item(2, 1).Value = "base_urn"
$ctrlSheet.Cells.Item(2, 2).Value = "urn:intuitem:risk:reference-controls:mitre-d3fend"

Write-Output "done"
